# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Swap the "Valor Mora" amounts between the first and last period rows
# of the statement table (Hoja1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 (period 2105): 44800 -> 56000
$ws.Range("F16").Value = 56000

# Row 22 (period 2011): 56000 -> 44800
$ws.Range("F22").Value = 44800
